# "add error ellipse parameters"
# The Points sheet tracked rmse_x / rmse_y / cov_xy for each point. This
# swaps the raw covariance column for the correlation coefficient
# (corr_xy) that is actually used to build the error ellipse, and fixes
# up the two sample values that went along with the old cov_xy column.

$wb = $excel.ActiveWorkbook

$points = $wb.Worksheets.Item("Points")
$points.Range("G1").Value = "corr_xy"
$points.Range("E3").Value = 2
$points.Range("F3").Value = 2

# -- Recreate the view/selection state left behind by the edit session --
$instruments = $wb.Worksheets.Item("Instruments")
$instruments.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$instruments.Range("F3").Select()

$measurements = $wb.Worksheets.Item("Measurements")
$measurements.Activate()
$measurements.Range("A4:XFD6").Select()

$points.Activate()
$points.Range("B27").Select()
